$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 13:22"

# Swap country labels for rows 11 (China) and 12 (Iran) -> after edit row 11 is Iran, row 12 is China
$ws.Range("A11").Value = "Iran"
$ws.Range("A12").Value = "China"

# Row 11 (Iran) - new figures
$ws.Range("B11").Value = 83505
$ws.Range("C11").Value = 1294
$ws.Range("D11").Value = 59273
$ws.Range("E11").Value = 19023
$ws.Range("F11").Value = 3389
$ws.Range("G11").Value = 91
$ws.Range("H11").Value = 5209

# Row 12 (China) - new figures
$ws.Range("B12").Value = 82747
$ws.Range("C12").Value = 12
$ws.Range("D12").Value = 77084
$ws.Range("E12").Value = 1031
$ws.Range("F12").Value = 81
$ws.Range("H12").Value = 4632

# Row 15 (Brasil) - new figures
$ws.Range("B15").Value = 39144
$ws.Range("C15").Value = 490
$ws.Range("E15").Value = 14530
$ws.Range("G15").Value = 22
$ws.Range("H15").Value = 2484

# Row 32 (Rumania) - new figures
$ws.Range("E32").Value = 6459
$ws.Range("G32").Value = 9
$ws.Range("H32").Value = 460

# Row 136 (Madagascar) - new figures
$ws.Range("D136").Value = 41
$ws.Range("E136").Value = 80

$wb.Save()
